# Remove the "Delta Period" (column I) and "Delta Payload Length" (column K)
# columns from the "Message Set" worksheet, shifting remaining columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Message Set")
$ws.Activate()

# Select the two (non-contiguous) whole columns that are being removed
# together (mimicking a ctrl-click multi-select of column I then column K),
# then delete them. Column K must be deleted first so that column I's index
# does not shift before it is removed.
$selection = $excel.Union($ws.Columns("I"), $ws.Columns("K"))
$selection.Select() | Out-Null

$ws.Columns("K").Delete()
$ws.Columns("I").Delete()

# Leave the final selection on column K (the last-activated column of the
# original multi-selection), matching the post-edit active cell.
$ws.Columns("K").Select() | Out-Null
